$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly rank report refresh (10th week) ---
# The sheet tracks one "rank event" column per week (B..H), newest week first.
# This week: drop the two oldest weekly columns (old B/Jun_24, old C/Jun_22),
# shift the remaining weekly columns left, and add two new weekly columns
# (Jun_27, Jun_26) at the front. Net effect: one fewer week column overall
# (7 -> 6), so the old trailing column (H) is removed at the end.

# 1) Shift row-1 week headers left by two (D..G <- old E..H)
for ($c = 4; $c -le 7; $c++) {
    $ws.Cells.Item(1, $c).Value = $ws.Cells.Item(1, $c + 1).Value()
}
# New two newest weeks at the front
$ws.Cells.Item(1, 2).Value = "Jun_27"
$ws.Cells.Item(1, 3).Value = "Jun_26"

# 2) Shift each stock's rank-event detail (column H) into the new last data
#    column (G) now that two columns have been dropped from the front.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value()
}

# 3) Drop the now-unused trailing column (H)
$ws.Range("H1").EntireColumn.Delete()

# 4) New group: add Benchmark / Evercore ISI rows
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
